# Absenzenliste-Template: widen the "Name" (surname) column by 2mm and
# narrow the "Vorname" (first name) column by 2mm.
#
# The document contains a single table whose grid column 3 ("Name") and
# column 4 ("Vorname") are both currently 1418 twips (70.9 pt) wide.
# They become 1548 twips (77.4 pt) and 1288 twips (64.4 pt) respectively.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Columns.Item(3).Width = 1548 / 20   # 1418 -> 1548 twips (Name / surname)
$t.Columns.Item(4).Width = 1288 / 20   # 1418 -> 1288 twips (Vorname / first name)
